# Add a new room-deletion record to the users table.
# Rows 28-32 stay blank (just re-assert the existing default formatting so the
# rows materialize in the sheet), and row 33 carries the new user/room record:
#   YAHYA | YAHYA | EA4C7814 | 1 | 105|103 | TRUE

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touch rows 28-32 so they exist with the same formatting the surrounding
# rows already use (column A/B/C/E/F = left aligned, column D = right
# aligned numeric), without introducing any new style entries.
$ws.Range("A28:C32").HorizontalAlignment = -4131
$ws.Range("E28:F32").HorizontalAlignment = -4131
$ws.Range("D28:D32").NumberFormat = "#,##0"

# Row 33 - new user/room entry.
$ws.Range("A33").Value2 = "YAHYA"
$ws.Range("B33").Value2 = "YAHYA"
$ws.Range("C33").Value2 = "EA4C7814"
$ws.Range("D33").Value2 = 1
$ws.Range("E33").Value2 = "105|103"

# F33 must hold the literal text "TRUE" (shared string), not a boolean.
# Assigning the string "TRUE" directly via Value/Value2 auto-coerces it to a
# Boolean cell, so build it as a text formula result in a scratch cell and
# paste-special the value in, which preserves the text type.
$ws.Range("H1").Formula = "=""TRUE"""
$ws.Range("H1").Copy()
$ws.Range("F33").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("H1").Clear()

Write-Host "Added rows 28-33"
